$d = $word.ActiveDocument

$d.Content.Find.Execute("38×38=1444", $true, $false, $false, $false, $false, $true, 1, $false, "38×16=608", 2)
$d.Content.Find.Execute("36×70=2520", $true, $false, $false, $false, $false, $true, 1, $false, "55×91=5005", 2)
$d.Content.Find.Execute("26×18=468", $true, $false, $false, $false, $false, $true, 1, $false, "24×61=1464", 2)
$d.Content.Find.Execute("99×25=2475", $true, $false, $false, $false, $false, $true, 1, $false, "65×45=2925", 2)
$d.Content.Find.Execute("89×69=6141", $true, $false, $false, $false, $false, $true, 1, $false, "62×82=5084", 2)
$d.Content.Find.Execute("13×31=403", $true, $false, $false, $false, $false, $true, 1, $false, "37×88=3256", 2)
$d.Content.Find.Execute("68×16=1088", $true, $false, $false, $false, $false, $true, 1, $false, "27×84=2268", 2)
$d.Content.Find.Execute("23×45=1035", $true, $false, $false, $false, $false, $true, 1, $false, "38×72=2736", 2)
$d.Content.Find.Execute("43×17=731", $true, $false, $false, $false, $false, $true, 1, $false, "71×80=5680", 2)
$d.Content.Find.Execute("12×31=372", $true, $false, $false, $false, $false, $true, 1, $false, "16×53=848", 2)
$d.Content.Find.Execute("61×85=5185", $true, $false, $false, $false, $false, $true, 1, $false, "69×23=1587", 2)
$d.Content.Find.Execute("70×61=4270", $true, $false, $false, $false, $false, $true, 1, $false, "91×30=2730", 2)
$d.Content.Find.Execute("39×90=3510", $true, $false, $false, $false, $false, $true, 1, $false, "64×68=4352", 2)
$d.Content.Find.Execute("48×45=2160", $true, $false, $false, $false, $false, $true, 1, $false, "48×18=864", 2)
$d.Content.Find.Execute("55×73=4015", $true, $false, $false, $false, $false, $true, 1, $false, "72×58=4176", 2)
$d.Content.Find.Execute("53×37=1961", $true, $false, $false, $false, $false, $true, 1, $false, "85×35=2975", 2)
$d.Content.Find.Execute("73×95=6935", $true, $false, $false, $false, $false, $true, 1, $false, "22×19=418", 2)
$d.Content.Find.Execute("21×80=1680", $true, $false, $false, $false, $false, $true, 1, $false, "36×50=1800", 2)
$d.Content.Find.Execute("35×25=875", $true, $false, $false, $false, $false, $true, 1, $false, "96×63=6048", 2)
$d.Content.Find.Execute("31×57=1767", $true, $false, $false, $false, $false, $true, 1, $false, "61×68=4148", 2)
$d.Content.Find.Execute("54×83=4482", $true, $false, $false, $false, $false, $true, 1, $false, "60×95=5700", 2)
$d.Content.Find.Execute("73×35=2555", $true, $false, $false, $false, $false, $true, 1, $false, "91×56=5096", 2)
$d.Content.Find.Execute("20×31=620", $true, $false, $false, $false, $false, $true, 1, $false, "60×70=4200", 2)
$d.Content.Find.Execute("78×23=1794", $true, $false, $false, $false, $false, $true, 1, $false, "71×14=994", 2)
$d.Content.Find.Execute("20×12=240", $true, $false, $false, $false, $false, $true, 1, $false, "35×56=1960", 2)
